# Quarterly financials update for DXF - "Doing Updates for Financials"
# Applies the refreshed quarterly figures to the Income Statement,
# Balance Sheet and Cash Flow Statement sections of the DXF sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DXF")

$ws.Range("D8").Value = 7500
$ws.Range("E8").Value = 17700
$ws.Range("J8").Value = 12300
$ws.Range("E9").Value = 5700
$ws.Range("J9").Value = 24700
$ws.Range("D10").Value = 5500
$ws.Range("E10").Value = 12100
$ws.Range("J10").Value = -12300
$ws.Range("D17").Value = 4800
$ws.Range("E17").Value = 10700
$ws.Range("J17").Value = 49100
$ws.Range("D18").Value = 2700
$ws.Range("E18").Value = 7000
$ws.Range("J18").Value = -36700
$ws.Range("E21").Value = 7200
$ws.Range("D23").Value = 2700
$ws.Range("E23").Value = 7200
$ws.Range("J23").Value = -36700
$ws.Range("E24").Value = 2200
$ws.Range("J24").Value = 7600
$ws.Range("E26").Value = 5000
$ws.Range("J26").Value = -44300
$ws.Range("E27").Value = 4000
$ws.Range("J27").Value = -44300
$ws.Range("F29").Value = -7100
$ws.Range("G29").Value = -50700
$ws.Range("H29").Value = -18900
$ws.Range("I29").Value = -88800
$ws.Range("E33").Value = 4000
$ws.Range("F33").Value = -7600
$ws.Range("G33").Value = -51200
$ws.Range("H33").Value = -19200
$ws.Range("I33").Value = -89100
$ws.Range("J33").Value = -44300
$ws.Range("E35").Value = 4000
$ws.Range("F35").Value = -7600
$ws.Range("G35").Value = -51200
$ws.Range("H35").Value = -19200
$ws.Range("I35").Value = -89100
$ws.Range("J35").Value = -44300
$ws.Range("H41").Value = 9400
$ws.Range("I41").Value = 30200
$ws.Range("D43").Value = 127200
$ws.Range("E43").Value = 123300
$ws.Range("H43").Value = 76500
$ws.Range("I43").Value = 96400
$ws.Range("H44").Value = 15800
$ws.Range("D45").Value = 4100
$ws.Range("E45").Value = 2300
$ws.Range("H45").Value = 6300
$ws.Range("I45").Value = 6800
$ws.Range("D46").Value = 131400
$ws.Range("E46").Value = 128800
$ws.Range("H46").Value = 108000
$ws.Range("I46").Value = 135600
$ws.Range("E48").Value = 100
$ws.Range("H48").Value = 500
$ws.Range("I49").Value = 2600
$ws.Range("F52").Value = 76400
$ws.Range("G52").Value = 81000
$ws.Range("D54").Value = 131500
$ws.Range("E54").Value = 129000
$ws.Range("F54").Value = 76400
$ws.Range("G54").Value = 81000
$ws.Range("H54").Value = 110300
$ws.Range("I54").Value = 138700
$ws.Range("H57").Value = 1900
$ws.Range("D58").Value = 34100
$ws.Range("E58").Value = 33600
$ws.Range("D59").Value = 5800
$ws.Range("E59").Value = 5700
$ws.Range("F59").Value = 26500
$ws.Range("G59").Value = 23400
$ws.Range("H59").Value = 18900
$ws.Range("I59").Value = 29100
$ws.Range("D60").Value = 39800
$ws.Range("E60").Value = 39300
$ws.Range("F60").Value = 26500
$ws.Range("G60").Value = 23400
$ws.Range("H60").Value = 20800
$ws.Range("I60").Value = 30100
$ws.Range("D66").Value = 58200
$ws.Range("E66").Value = 57300
$ws.Range("F66").Value = 26500
$ws.Range("G66").Value = 23400
$ws.Range("H66").Value = 20800
$ws.Range("I66").Value = 30100
$ws.Range("D72").Value = 15700
$ws.Range("E72").Value = 14000
$ws.Range("F72").Value = -24600
$ws.Range("G72").Value = -16900
$ws.Range("H72").Value = 15100
$ws.Range("I72").Value = 34300
$ws.Range("D76").Value = 73400
$ws.Range("E76").Value = 71700
$ws.Range("F76").Value = 49900
$ws.Range("G76").Value = 57600
$ws.Range("H76").Value = 89500
$ws.Range("I76").Value = 108600
$ws.Range("E81").Value = 4000
$ws.Range("F81").Value = -7600
$ws.Range("G81").Value = -51200
$ws.Range("H81").Value = -19200
$ws.Range("I81").Value = -89100
$ws.Range("J81").Value = -44300
$ws.Range("D89").Value = 1700
$ws.Range("E89").Value = 6600
$ws.Range("G89").Value = -5100
$ws.Range("H89").Value = -23300
$ws.Range("I89").Value = -107300
$ws.Range("D91").Value = 0
$ws.Range("D94").Value = -5300
$ws.Range("E94").Value = -20200
$ws.Range("G94").Value = -19900
$ws.Range("H94").Value = 2500
$ws.Range("I94").Value = -8700
$ws.Range("E100").Value = 2500
$ws.Range("E102").Value = -11100
$ws.Range("F102").Value = 2600
$ws.Range("G102").Value = -25000
$ws.Range("H102").Value = -20800
$ws.Range("I102").Value = -116000
